$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column H (8th column), shifting H:EU -> I:EV
$ws.Columns.Item(8).Insert()

# New header for the inserted column
$ws.Cells.Item(1, 8).Value = "Hire Date"

# Restore selection to H2, matching the sheetView's active cell after the edit
$ws.Range("H2").Select() | Out-Null
